# Reorder the (Title, Abstract) pairs in rows 2-5 so that each title keeps
# its own abstract, but the pairs cyclically shift up one row (row3->row2,
# row4->row3, row5->row4, row2->row5), and set the Score in D2 to 0.
#
# Note: this runtime's Range.Value getter doesn't return the underlying
# variant correctly, so Value2 is used for reads (writes use Value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$titleA2 = $ws.Range("A2").Value2
$titleA3 = $ws.Range("A3").Value2
$titleA4 = $ws.Range("A4").Value2
$titleA5 = $ws.Range("A5").Value2

$abstractB2 = $ws.Range("B2").Value2
$abstractB3 = $ws.Range("B3").Value2
$abstractB4 = $ws.Range("B4").Value2
$abstractB5 = $ws.Range("B5").Value2

# New row 2 gets old row 3's title+abstract, new row 3 gets old row 4's,
# new row 4 gets old row 5's, and new row 5 gets old row 2's.
$ws.Range("A2").Value = $titleA3
$ws.Range("B2").Value = $abstractB3

$ws.Range("A3").Value = $titleA4
$ws.Range("B3").Value = $abstractB4

$ws.Range("A4").Value = $titleA5
$ws.Range("B4").Value = $abstractB5

$ws.Range("A5").Value = $titleA2
$ws.Range("B5").Value = $abstractB2

$ws.Range("D2").Value = 0

$wb.Save()
